# Updated cryptos list on Thu Jun 22 18:51:17 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.038.25"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("D4").Value = "'0.9985"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "'244.19"
$ws.Range("E5").Value = "  -2.26%  "

$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").Value = "'0.4960"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").Value = "'44.34"
$ws.Range("E8").Value = "  -2.93%  "

$ws.Range("D9").Value = "'0.2921"
$ws.Range("E9").Value = "  +2.52%  "

$ws.Range("D10").Value = "'0.06622"
$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D11").Value = "1.879.58"
$ws.Range("E11").Value = "  +0.57%  "

$ws.Range("D12").Value = "'16.94"
$ws.Range("E12").Value = "  -0.73%  "

$ws.Range("D13").Value = "'0.07200"
$ws.Range("E13").Value = "  -0.01%  "

$ws.Range("D14").Value = "'0.6674"
$ws.Range("E14").Value = "  +0.81%  "

$ws.Range("D15").Value = "'86.02"
$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").Value = "'4.844"
$ws.Range("E16").Value = "  +0.61%  "

$ws.Range("D17").Value = "30.011.55"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "'0.000007821"
$ws.Range("E18").Value = "  +3.51%  "

$ws.Range("D19").Value = "'0.9981"
$ws.Range("E19").Value = "  -0.02%  "

$ws.Range("D20").Value = "'12.82"
$ws.Range("E20").Value = "  -0.94%  "

$ws.Range("D21").Value = "2.121.52"
$ws.Range("E21").Value = "  +0.58%  "

$ws.Range("D22").Value = "'0.9982"
$ws.Range("E22").Value = "  -0.24%  "

$ws.Range("D23").Value = "'4.767"
$ws.Range("E23").Value = "  -0.19%  "

# Row 24
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'5.602"
$ws.Range("E24").Value = "  +1.74%  "

# Row 25
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'9.163"
$ws.Range("E25").Value = "  +1.18%  "

$ws.Range("D26").Value = "'150.09"
$ws.Range("E26").Value = "  +3.69%  "

$ws.Range("D27").Value = "'135.85"
$ws.Range("E27").Value = "  -0.09%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("D29").Value = "'1.907"
$ws.Range("E29").Value = "  -2.42%  "

$ws.Range("E30").Value = "  -2.05%  "

$ws.Range("D31").Value = "'4.181"
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("D32").Value = "'0.08673"
$ws.Range("E32").Value = "  +0.52%  "

$ws.Range("D33").Value = "'3.955"
$ws.Range("E33").Value = "  +1.31%  "

$ws.Range("D34").Value = "'0.04987"
$ws.Range("E34").Value = "  -1.74%  "

$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("D36").Value = "'0.7023"
$ws.Range("E36").Value = "  +2.11%  "

$ws.Range("D37").Value = "'2.652"
$ws.Range("E37").Value = "  -1.48%  "

$ws.Range("D38").Value = "'2.200"
$ws.Range("E38").Value = "  -6.37%  "

$ws.Range("E39").Value = "  -1.97%  "

$ws.Range("D40").Value = "'0.9352"
$ws.Range("E40").Value = "  -2.53%  "

$ws.Range("D41").Value = "'0.01642"
$ws.Range("E41").Value = "  +0.69%  "

$ws.Range("D42").Value = "'5.966"
$ws.Range("E42").Value = "  -2.25%  "

$ws.Range("D43").Value = "'0.9995"
$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("D44").Value = "'0.4187"
$ws.Range("E44").Value = "  -0.21%  "

$ws.Range("D45").Value = "'101.52"
$ws.Range("E45").Value = "  -1.85%  "

$ws.Range("D46").Value = "'7.540"
$ws.Range("E46").Value = "  +0.62%  "

$ws.Range("D47").Value = "'0.1261"
$ws.Range("E47").Value = "  +0.37%  "

$ws.Range("D48").Value = "'0.05722"
$ws.Range("E48").Value = "  +1.44%  "

$ws.Range("D49").Value = "'32.36"
$ws.Range("E49").Value = "  -0.73%  "

$ws.Range("D50").Value = "'8.216"
$ws.Range("E50").Value = "  -1.18%  "

$ws.Range("D51").Value = "'55.85"
$ws.Range("E51").Value = "  +1.71%  "
